# Management command: generate Excel dummy transactions
# Regenerates the purchase-request seed data with a larger, more varied
# dummy dataset (10 rows instead of 2) and tightens the CODE/date display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dummy "purchase request" transactions: CODE, DATE (serial), REQUESTOR,
# APPROVER, VENDOR, STATUS
$transactions = @(
    @("P-REQ-1",  45658, "TEST-ID-15", "TEST-ID-15", "VENDOR-02", "CLOSED"),
    @("P-REQ-2",  45659, "TEST-ID-15", "TEST-ID-15", "VENDOR-02", "CLOSED"),
    @("P-REQ-3",  45660, "TEST-ID-14", "TEST-ID-15", "VENDOR-01", "CLOSED"),
    @("P-REQ-4",  45661, "TEST-ID-14", "TEST-ID-14", "VENDOR-02", "CLOSED"),
    @("P-REQ-5",  45667, "TEST-ID-14", "TEST-ID-14", "VENDOR-02", "CLOSED"),
    @("P-REQ-6",  45671, "TEST-ID-15", "TEST-ID-14", "VENDOR-02", "CLOSED"),
    @("P-REQ-7",  45673, "TEST-ID-15", "TEST-ID-14", "VENDOR-02", "CLOSED"),
    @("P-REQ-8",  45679, "TEST-ID-15", "TEST-ID-15", "VENDOR-01", "CLOSED"),
    @("P-REQ-9",  45685, "TEST-ID-14", "TEST-ID-14", "VENDOR-02", "CLOSED"),
    @("P-REQ-10", 45688, "TEST-ID-14", "TEST-ID-15", "VENDOR-01", "CLOSED")
)

$headerRow = 1
$row = $headerRow + 1
foreach ($tx in $transactions) {
    $ws.Cells.Item($row, 1).Value = $tx[0]
    $ws.Cells.Item($row, 2).Value = $tx[1]
    $ws.Cells.Item($row, 3).Value = $tx[2]
    $ws.Cells.Item($row, 4).Value = $tx[3]
    $ws.Cells.Item($row, 5).Value = $tx[4]
    $ws.Cells.Item($row, 6).Value = $tx[5]
    $row++
}

$lastRow = $headerRow + $transactions.Count

# Plain date display (drop the old ";@" text fallback segment)
$ws.Range("B" + ($headerRow + 1) + ":B" + $lastRow).NumberFormat = "yyyy\-mm\-dd"

# Leave the selection where the new STATUS column was just filled in
$null = $ws.Range("F" + ($headerRow + 1) + ":F" + $lastRow).Select()
